$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a thin border (all four sides) around the whole table range,
# fixing up the table's border as noted in the commit message.
$range = $ws.Range("A1:E14")
$range.Borders.LineStyle = 1      # xlContinuous
$range.Borders.Weight = 2         # xlThin

# Update the active selection to match the edited workbook's saved view state.
$ws.Range("G10").Select()
